# The source diff for this revision only reorders XML attributes/namespace
# declarations (e.g. alphabetizing w:pgSz, w:pgMar, w:rFonts, w:lang,
# w:latentStyles/w:lsdException and w:style attributes, plus the root
# <w:document> namespace declarations) as a side effect of the Apache POI
# packaging/version upgrade mentioned in the commit message. Every
# attribute value and every piece of document content is identical before
# and after the change - there is no semantic edit to replicate through
# the Word object model. Touching the document here would only risk
# introducing unintended differences, so we deliberately leave the
# content untouched.
$d = $word.ActiveDocument
